$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 829.9091
$ws.Range("I2").Value = 282.7143
$ws.Range("K2").Value = 282.7143
$ws.Range("M2").Value = -169.7143

$ws.Range("H12").Value = 483.25
$ws.Range("I12").Value = 749.5
$ws.Range("K12").Value = 749.5
$ws.Range("M12").Value = -579.5

$ws.Range("H17").Value = 206451.67
$ws.Range("J17").Value = 210718.11
$ws.Range("L17").Value = 632154.33
$ws.Range("N17").Value = -632490.33

$ws.Range("H58").Value = 94.111115
$ws.Range("I58").Value = 94.111115
$ws.Range("K58").Value = 282.333345
$ws.Range("M58").Value = -132.333345

$ws.Range("H92").Value = 1289.4286
$ws.Range("I92").Value = 996.44446
$ws.Range("K92").Value = 996.44446
$ws.Range("M92").Value = 251.55554

$ws.Range("H100").Value = 8387.958000000001
$ws.Range("I100").Value = 1543.2222
$ws.Range("J100").Value = 12494.8
$ws.Range("K100").Value = 1543.2222
$ws.Range("L100").Value = 12494.8
$ws.Range("M100").Value = -1002.2222
$ws.Range("N100").Value = -13576.8

$ws.Range("H107").Value = 573.56525
$ws.Range("I107").Value = 713.58826
$ws.Range("J107").Value = 176.83333
$ws.Range("K107").Value = 713.58826
$ws.Range("L107").Value = 176.83333
$ws.Range("M107").Value = 1206.41174
$ws.Range("N107").Value = -4016.83333

$ws.Range("H113").Value = 4438.3335
$ws.Range("J113").Value = 5499.6
$ws.Range("L113").Value = 5499.6
$ws.Range("N113").Value = -12007.6

$ws.Range("H132").Value = 11483.053
$ws.Range("I132").Value = 2330.2766
$ws.Range("K132").Value = 6990.8298
$ws.Range("M132").Value = -4460.8298

$ws.Range("H138").Value = 168335.5
$ws.Range("I138").Value = 1002.15
$ws.Range("J138").Value = 252002.17
$ws.Range("K138").Value = 3006.45
$ws.Range("L138").Value = 756006.51
$ws.Range("M138").Value = 2133.55
$ws.Range("N138").Value = -766286.51

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5636.469
$ws.Range("I32").Value = 5729.0947
$ws.Range("K32").Value = 5729.0947
$ws.Range("M32").Value = -5442.0947

$ws.Range("H97").Value = 1157.0476
$ws.Range("I97").Value = 1164.95
$ws.Range("K97").Value = 1164.95
$ws.Range("M97").Value = -668.95

$ws.Range("H122").Value = 2909.1702
$ws.Range("J122").Value = 4021.2727
$ws.Range("L122").Value = 12063.8181
$ws.Range("N122").Value = -16963.8181

$ws.Range("H132").Value = 2273.5476
$ws.Range("I132").Value = 2217.0132
$ws.Range("J132").Value = 2810.625
$ws.Range("K132").Value = 6651.0396
$ws.Range("L132").Value = 8431.875
$ws.Range("M132").Value = -4121.0396
$ws.Range("N132").Value = -13491.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2801.85
$ws.Range("I99").Value = 2135
$ws.Range("K99").Value = 2135
$ws.Range("M99").Value = -637

$ws.Range("H105").Value = 3521.5312
$ws.Range("I105").Value = 3378
$ws.Range("J105").Value = 3888.3333
$ws.Range("K105").Value = 3378
$ws.Range("L105").Value = 3888.3333
$ws.Range("M105").Value = -1631
$ws.Range("N105").Value = -7382.3333

$ws.Range("H107").Value = 1264.8695
$ws.Range("I107").Value = 864.8
$ws.Range("K107").Value = 864.8
$ws.Range("M107").Value = 1055.2

$ws.Range("H126").Value = 100000
$ws.Range("J126").Value = 100000
$ws.Range("L126").Value = 100000
$ws.Range("N126").Value = -109880

$ws.Range("H134").Value = 6762.9443
$ws.Range("I134").Value = 3497.5715
$ws.Range("K134").Value = 10492.7145
$ws.Range("M134").Value = -7957.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 684
$ws.Range("I22").Value = 674.4
$ws.Range("K22").Value = 674.4
$ws.Range("M22").Value = -324.4

$ws.Range("H58").Value = 1813.85
$ws.Range("I58").Value = 1828.7333
$ws.Range("K58").Value = 1828.7333
$ws.Range("M58").Value = -1625.7333

$ws.Range("H132").Value = 1820878
$ws.Range("I132").Value = 2502696.8
$ws.Range("J132").Value = 2694.5
$ws.Range("K132").Value = 7508090.399999999
$ws.Range("L132").Value = 8083.5
$ws.Range("M132").Value = -7505560.399999999
$ws.Range("N132").Value = -13143.5

$ws.Range("H134").Value = 4662.4443
$ws.Range("I134").Value = 1936
$ws.Range("K134").Value = 5808
$ws.Range("M134").Value = -3273

$ws.Range("H136").Value = 1813.85
$ws.Range("I136").Value = 1828.7333
$ws.Range("K136").Value = 5486.199900000001
$ws.Range("M136").Value = -2936.199900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 173.14285
$ws.Range("J2").Value = 135.71428
$ws.Range("L2").Value = 814.28568
$ws.Range("N2").Value = -1040.28568

$ws.Range("H38").Value = 4515.6924
$ws.Range("J38").Value = 13944.25
$ws.Range("L38").Value = 41832.75
$ws.Range("N38").Value = -42526.75

$ws.Range("H103").Value = 813.8
$ws.Range("J103").Value = 1085.3334
$ws.Range("L103").Value = 3256.0002
$ws.Range("N103").Value = -5014.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 24999
$ws.Range("J26").Value = 24999
$ws.Range("L26").Value = 24999
$ws.Range("N26").Value = -25559

$ws.Range("H50").Value = 24999
$ws.Range("J50").Value = 24999
$ws.Range("L50").Value = 24999
$ws.Range("N50").Value = -25995

$ws.Range("H122").Value = 2131.3635
$ws.Range("I122").Value = 1642.875
$ws.Range("J122").Value = 3434
$ws.Range("K122").Value = 4928.625
$ws.Range("L122").Value = 10302
$ws.Range("M122").Value = -2478.625
$ws.Range("N122").Value = -15202

$ws.Range("H132").Value = 4764432
$ws.Range("I132").Value = 5850368
$ws.Range("J132").Value = 3021
$ws.Range("K132").Value = 17551104
$ws.Range("L132").Value = 9063
$ws.Range("M132").Value = -17548574
$ws.Range("N132").Value = -14123

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8477.951999999999
$ws.Range("I46").Value = 4145.1665
$ws.Range("J46").Value = 10211.066
$ws.Range("K46").Value = 4145.1665
$ws.Range("L46").Value = 10211.066
$ws.Range("M46").Value = -3957.1665
$ws.Range("N46").Value = -10587.066

$ws.Range("H93").Value = 2057.6667
$ws.Range("J93").Value = 2099
$ws.Range("L93").Value = 2099
$ws.Range("N93").Value = -4595

$ws.Range("H132").Value = 2626.0159
$ws.Range("I132").Value = 2552.6
$ws.Range("J132").Value = 2908.3845
$ws.Range("K132").Value = 7657.799999999999
$ws.Range("L132").Value = 8725.1535
$ws.Range("M132").Value = -5127.799999999999
$ws.Range("N132").Value = -13785.1535

$ws.Range("H136").Value = 4034.2144
$ws.Range("I136").Value = 3249.8333
$ws.Range("K136").Value = 9749.499899999999
$ws.Range("M136").Value = -7199.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 18312
$ws.Range("J45").Value = 18312
$ws.Range("L45").Value = 18312
$ws.Range("N45").Value = -19294

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H109").Value = 87499.60000000001
$ws.Range("J109").Value = 87499.60000000001
$ws.Range("L109").Value = 87499.60000000001
$ws.Range("N109").Value = -90273.60000000001
